$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three obsolete rows (old rows 5-7: EC/FAP/MuSC x MuSC target
# combinations that no longer exist after the TPM recalculation).
$ws.Rows("5:7").Delete()

# Row 2: ECs -> MuSCs (Target cluster changed from ECs to MuSCs), updated TPM values
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("M2").Value = 0.5233716666666667
$ws.Range("N2").Value = 1.570115
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 33.44930201643389
$ws.Range("R2").Value = 301.0437181479049
$ws.Range("S2").Value = 0.4067926910433548
$ws.Range("T2").Value = 0.4067926910433549

# Row 3: Sending cluster changes FAPs (was ECs), target stays MuSCs, updated TPM values
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 57.4434
$ws.Range("H3").Value = 172.3302
$ws.Range("I3").Value = 0.3656254573230189
$ws.Range("J3").Value = 0.365625457323019
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5233716666666667
$ws.Range("N3").Value = 1.570115
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 30.064247997
$ws.Range("R3").Value = 270.578231973
$ws.Range("S3").Value = 0.3656254573230189
$ws.Range("T3").Value = 0.365625457323019

# Row 4: Sending cluster changes to MuSCs (was FAPs), target cluster changes to MuSCs (was ECs)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 35.755375
$ws.Range("H4").Value = 107.266125
$ws.Range("I4").Value = 0.2275818516336261
$ws.Range("J4").Value = 0.2275818516336262
$ws.Range("M4").Value = 0.5233716666666667
$ws.Range("N4").Value = 1.570115
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 18.71335020604167
$ws.Range("R4").Value = 168.420151854375
$ws.Range("S4").Value = 0.2275818516336261
$ws.Range("T4").Value = 0.2275818516336262
